$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1113.3704
$ws.Range("I19").Value = 188.25
$ws.Range("K19").Value = 188.25
$ws.Range("M19").Value = -13.25
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = ""
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7800.5
$ws.Range("I61").Value = 8500.75
$ws.Range("K61").Value = 8500.75
$ws.Range("M61").Value = -8288.75
$ws.Range("H74").Value = 2214.1538
$ws.Range("I74").Value = 1421.6666
$ws.Range("K74").Value = 1421.6666
$ws.Range("M74").Value = -547.6666
$ws.Range("H77").Value = 2214.1538
$ws.Range("I77").Value = 1421.6666
$ws.Range("K77").Value = 7108.333000000001
$ws.Range("M77").Value = -2740.333000000001
$ws.Range("H122").Value = 12076.174
$ws.Range("I122").Value = 8187.6
$ws.Range("K122").Value = 24562.8
$ws.Range("M122").Value = -22112.8
$ws.Range("H132").Value = 5760.615
$ws.Range("I132").Value = 4756.4287
$ws.Range("J132").Value = 6932.1665
$ws.Range("K132").Value = 14269.2861
$ws.Range("L132").Value = 20796.4995
$ws.Range("M132").Value = -11739.2861
$ws.Range("N132").Value = -25856.4995
$ws.Range("H136").Value = 7800.5
$ws.Range("I136").Value = 8500.75
$ws.Range("K136").Value = 25502.25
$ws.Range("M136").Value = -22952.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19708
$ws.Range("H134").Value = 2833.75
$ws.Range("I134").Value = 2326.4614
$ws.Range("K134").Value = 6979.3842
$ws.Range("M134").Value = -4444.3842
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 123.85714
$ws.Range("J7").Value = 362
$ws.Range("K7").Value = 123.85714
$ws.Range("L7").Value = 362
$ws.Range("M7").Value = -10.85714
$ws.Range("N7").Value = -588
$ws.Range("H12").Value = 750
$ws.Range("I12").Value = 750
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -580
$ws.Range("N12").Value = ""
$ws.Range("H31").Value = 2136.9333
$ws.Range("I31").Value = 1466.2693
$ws.Range("K31").Value = 1466.2693
$ws.Range("M31").Value = -1171.2693
$ws.Range("H34").Value = 2136.9333
$ws.Range("I34").Value = 1466.2693
$ws.Range("K34").Value = 1466.2693
$ws.Range("M34").Value = -1264.2693
$ws.Range("H42").Value = 9006.799999999999
$ws.Range("I42").Value = 4600
$ws.Range("J42").Value = 15617
$ws.Range("K42").Value = 4600
$ws.Range("L42").Value = 15617
$ws.Range("M42").Value = -4007
$ws.Range("N42").Value = -16803
$ws.Range("H58").Value = 3113.818
$ws.Range("I58").Value = 1958.1428
$ws.Range("K58").Value = 1958.1428
$ws.Range("M58").Value = -1755.1428
$ws.Range("H99").Value = 2282.6155
$ws.Range("I99").Value = 2282.6155
$ws.Range("K99").Value = 2282.6155
$ws.Range("M99").Value = -784.6154999999999
$ws.Range("H126").Value = 2282.6155
$ws.Range("I126").Value = 2282.6155
$ws.Range("K126").Value = 6847.8465
$ws.Range("M126").Value = -4377.8465
$ws.Range("H132").Value = 1805.0952
$ws.Range("I132").Value = 1732
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 5196
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -2666
$ws.Range("N132").Value = -12558.5
$ws.Range("H134").Value = 2880.1765
$ws.Range("I134").Value = 2536.1538
$ws.Range("K134").Value = 7608.4614
$ws.Range("M134").Value = -5073.4614
$ws.Range("H136").Value = 3113.818
$ws.Range("I136").Value = 1958.1428
$ws.Range("K136").Value = 5874.428400000001
$ws.Range("M136").Value = -3324.428400000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10001
$ws.Range("I57").Value = 10001
$ws.Range("K57").Value = 30003
$ws.Range("M57").Value = -29444
$ws.Range("H137").Value = 6006.6
$ws.Range("I137").Value = 7500
$ws.Range("J137").Value = 5011
$ws.Range("K137").Value = 22500
$ws.Range("L137").Value = 15033
$ws.Range("M137").Value = -17400
$ws.Range("N137").Value = -25233
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 236.41667
$ws.Range("I2").Value = 133.5
$ws.Range("J2").Value = 287.875
$ws.Range("K2").Value = 133.5
$ws.Range("L2").Value = 287.875
$ws.Range("M2").Value = -20.5
$ws.Range("N2").Value = -513.875
$ws.Range("H70").Value = 24996.334
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("H73").Value = 24996.334
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("H132").Value = 3263.4546
$ws.Range("I132").Value = 3294.5264
$ws.Range("K132").Value = 9883.5792
$ws.Range("M132").Value = -7353.5792
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1464.3334
$ws.Range("I16").Value = 1595.5
$ws.Range("J16").Value = 1202
$ws.Range("K16").Value = 1595.5
$ws.Range("L16").Value = 1202
$ws.Range("M16").Value = -1425.5
$ws.Range("N16").Value = -1542
$ws.Range("H46").Value = 4100
$ws.Range("J46").Value = 6142.857
$ws.Range("L46").Value = 6142.857
$ws.Range("N46").Value = -6518.857
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("H63").Value = 30249
$ws.Range("J63").Value = 30249
$ws.Range("L63").Value = 30249
$ws.Range("N63").Value = -31497
$ws.Range("H65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H66").Value = 30249
$ws.Range("J66").Value = 30249
$ws.Range("L66").Value = 90747
$ws.Range("N66").Value = -96987
$ws.Range("H105").Value = 17450
$ws.Range("J105").Value = 17450
$ws.Range("L105").Value = 17450
$ws.Range("N105").Value = -24438
